$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.257.85'
$ws.Range('E2').Value = '  +0.52%  '

$ws.Range('D3').Value = '1.590.41'
$ws.Range('E3').Value = '  +1.20%  '

$ws.Range('E4').Value = '  -0.16%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '212.29'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.74%  '

$ws.Range('E6').Value = '  +0.80%  '

$ws.Range('E7').Value = '  -0.19%  '

$ws.Range('E8').Value = '  +0.80%  '

$ws.Range('E9').Value = '  +0.09%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.39'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.04%  '

$ws.Range('E11').Value = '  +0.54%  '

$ws.Range('D12').Value = '1.813.39'
$ws.Range('E12').Value = '  +1.17%  '

$ws.Range('D13').Value = '1.593.59'
$ws.Range('E13').Value = '  +2.01%  '

$ws.Range('E14').Value = '  -0.02%  '

$ws.Range('E15').Value = '  +1.74%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.39'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.43%  '

$ws.Range('D17').Value = '26.254.61'
$ws.Range('E17').Value = '  +0.52%  '

$ws.Range('E18').Value = '  +0.32%  '

$ws.Range('E19').Value = '  +2.16%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '213.72'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +3.26%  '

$ws.Range('E21').Value = '  -0.13%  '

$ws.Range('E22').Value = '  +1.26%  '

$ws.Range('E23').Value = '  +2.08%  '

$ws.Range('E24').Value = '  -2.14%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '143.72'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.03%  '

$ws.Range('E26').Value = '  -0.17%  '

$ws.Range('E27').Value = '  +1.40%  '

$ws.Range('E28').Value = '  -0.01%  '

$ws.Range('E29').Value = '  +0.14%  '

$ws.Range('E30').Value = '  -0.90%  '

$ws.Range('E31').Value = '  +1.46%  '

$ws.Range('E32').Value = '  +0.13%  '

$ws.Range('D33').Value = '1.340.42'
$ws.Range('E33').Value = '  +4.84%  '

$ws.Range('E34').Value = '  -0.90%  '

$ws.Range('E35').Value = '  +0.04%  '

$ws.Range('E36').Value = '  +0.26%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.591'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.77%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0166'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.67%  '

$ws.Range('E39').Value = '  +0.63%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.77'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +3.70%  '

$ws.Range('E41').Value = '  -0.11%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.01'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -7.54%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.14'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.53%  '

$ws.Range('E44').Value = '  +0.73%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '61.95'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.21%  '

$ws.Range('D46').Value = '1.725.07'
$ws.Range('E46').Value = '  +1.15%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '85.80'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.60%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.47'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.80%  '

$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0502'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.50%  '

$ws.Range('B50').Value = 'Algorand'
$ws.Range('C50').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0976'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.37%  '

$ws.Range('E51').Value = '  -0.30%  '
